$wb = $excel.ActiveWorkbook

# --- Rename the first sheet (org_molgenis_test_TypeTest -> TypeTest) ---
$wsTypeTest = $wb.Worksheets.Item("org_molgenis_test_TypeTest")
$wsTypeTest.Name = "TypeTest"

# --- entities sheet: drop the stray B2 cell, move selection to B7 ---
$wsEntities = $wb.Worksheets.Item("entities")
$wsEntities.Range("B2").ClearContents()

# --- attributes sheet: every "entity" cell in column B (rows 4-50) pointed
#     at the old fully-qualified entity name; repoint it at the new short
#     sheet/entity name "TypeTest" ---
$wsAttributes = $wb.Worksheets.Item("attributes")
for ($r = 4; $r -le 50; $r++) {
    $wsAttributes.Cells.Item($r, 2).Value = "TypeTest"
}

# Reset the font color (black) on the "defaultValue" column for the boolean
# rows that got the style bump (D36, D37, D42, D43, D48)
$boolDefaultRows = 36, 37, 42, 43, 48
foreach ($r in $boolDefaultRows) {
    $wsAttributes.Cells.Item($r, 4).Font.Color = 0
}

# --- Update selections / active sheet to match the saved view state ---
$wsTypeTestRef = $wb.Worksheets.Item("org_molgenis_test_TypeTestRef")
$wsTypeTestRef.Activate()
$wsTypeTestRef.Range("F49").Select()

$wsEntities.Activate()
$wsEntities.Range("B7").Select()

$wsPackages = $wb.Worksheets.Item("packages")
$wsPackages.Activate()
$wsPackages.Range("A2").Select()

$wsAttributes.Activate()
$wsAttributes.Range("B3").Select()
